$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New vm_pu values for rows 2-25 (bus results), columns B-F and I-N.
# Column G (and H, which is blank) are unchanged.
$columns = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
    2 = @(1.02, 1.016781629710088, 1.04640570893408, 1.0182291542689, 1.050219672561107, 1.039878594663416, 1.021999482056776, 1.049171235683644, 1.02107599100649, 1.052974556895496, 1.011315331217632)
    3 = @(1.02, 1.018273729989685, 1.047100559540913, 1.019510530630493, 1.05119956134051, 1.040064036485619, 1.023124330810712, 1.049678025853052, 1.022161822076517, 1.053766416524224, 1.011697610754436)
    4 = @(1.02, 1.019237761502197, 1.047547989523303, 1.020338709891954, 1.051831198254204, 1.040181279457985, 1.023850410904467, 1.050003097092457, 1.022862936438564, 1.054275775879018, 1.011944086716067)
    5 = @(1.02, 1.019642698901204, 1.047735565734536, 1.020686652641428, 1.052096160871901, 1.040229909726947, 1.024155235782094, 1.050139073206841, 1.023157332763208, 1.054489186778733, 1.012047495136177)
    6 = @(1.02, 1.019710669863488, 1.047767029906365, 1.020745060788263, 1.052140615422295, 1.040238036330069, 1.024206392785732, 1.050161864089723, 1.023206742689343, 1.05452497698887, 1.012064845594266)
    7 = @(1.02, 1.019243173628108, 1.047550497982761, 1.020343359992948, 1.051834740965928, 1.04018193184657, 1.023854485630316, 1.050004916699258, 1.022866871554523, 1.054278630328966, 1.011945469287842)
    8 = @(1.02, 1.017286198161456, 1.046640989165053, 1.018662402736702, 1.05055133157019, 1.039941834917569, 1.0223799999954, 1.049343099061025, 1.021443264300943, 1.053242796802427, 1.011444708518428)
    9 = @(1.02, 1.013826264588717, 1.045021595611019, 1.015692774431264, 1.048271242920713, 1.039497698904769, 1.019767964847435, 1.048155021936845, 1.018923051867951, 1.051394303177267, 1.010555454847119)
    10 = @(1.02, 1.011511446996195, 1.04393076380479, 1.013707599035333, 1.046738631390412, 1.039187461266105, 1.018017003495174, 1.047348263961936, 1.017234783392103, 1.050146295806834, 1.009957907501216)
    11 = @(1.02, 1.01050705665331, 1.043455756055927, 1.012846634293134, 1.04607199735782, 1.039049770195881, 1.01725646676491, 1.046995438176097, 1.016501748339656, 1.049602159592727, 1.009698021301186)
    12 = @(1.02, 1.010133662872713, 1.043278915581638, 1.012526621115517, 1.045823926742409, 1.03899812130388, 1.016973608760918, 1.046863857584019, 1.016229159449954, 1.04939947951286, 1.009601313963183)
    13 = @(1.02, 1.010213771665905, 1.043316866634549, 1.012595274752927, 1.045877159263782, 1.039009222988632, 1.017034299193293, 1.046892105849604, 1.016287644683016, 1.049442980612042, 1.009622065925039)
    14 = @(1.02, 1.010476198357099, 1.043441146558191, 1.012820186278186, 1.046051501016558, 1.039045511170438, 1.017233093031504, 1.04698457240358, 1.016479222355317, 1.049585417502412, 1.009690031004996)
    15 = @(1.02, 1.01063784562542, 1.0435176663169, 1.012958733329901, 1.046158858633459, 1.039067802689712, 1.017355528474307, 1.047041474467393, 1.016597218816628, 1.049673102830412, 1.009731883409186)
    16 = @(1.02, 1.011578060915385, 1.043962232222477, 1.01376470898827, 1.046782810265198, 1.039196528653276, 1.018067427568801, 1.047371606193176, 1.017283389788193, 1.050182329381635, 1.009975131003732)
    17 = @(1.02, 1.012167275433644, 1.044240381570786, 1.014269904605783, 1.047173393317201, 1.039276376473276, 1.018513347016154, 1.047577753313696, 1.017713266050916, 1.05050075108138, 1.010127405952957)
    18 = @(1.02, 1.012510756246865, 1.044402363748355, 1.014564444860141, 1.047400924176416, 1.039322626546065, 1.018773217121003, 1.047697658327908, 1.017963812686997, 1.050686120256254, 1.010216115095047)
    19 = @(1.02, 1.012627841010677, 1.044457551792439, 1.014664853317975, 1.047478457248658, 1.039338341703526, 1.018861787772465, 1.047738485670091, 1.018049210073996, 1.050749265202418, 1.010246343981831)
    20 = @(1.02, 1.012104078843182, 1.044210565428016, 1.01421571557795, 1.047131517423245, 1.039267843051606, 1.018465527622654, 1.047555670550883, 1.017667164393868, 1.050466624779354, 1.010111079712053)
    21 = @(1.02, 1.010398929102922, 1.04340456030244, 1.012753961360213, 1.046000174237528, 1.039034839127213, 1.017174563212489, 1.046957357814865, 1.016422816029889, 1.049543488961737, 1.009670021804197)
    22 = @(1.02, 1.009324986534748, 1.042895470081477, 1.01183366547138, 1.045286231673529, 1.038885422540334, 1.016360791960544, 1.04657813489223, 1.01563866541253, 1.048959814796082, 1.009391703231763)
    23 = @(1.02, 1.009894481978082, 1.043165568712072, 1.012321650829092, 1.045664955328275, 1.038964907704483, 1.016792387823412, 1.046779456441268, 1.016054529178212, 1.049269541284565, 1.009539341429055)
    24 = @(1.02, 1.012132635274322, 1.044224038856053, 1.014240201674623, 1.047150440232652, 1.039271699938314, 1.018487135850375, 1.047565649842901, 1.017687996344593, 1.050482046101649, 1.010118457178355)
    25 = @(1.02, 1.014722150259201, 1.045442228449873, 1.016461425609722, 1.048862907176831, 1.039615012749204, 1.020444907460198, 1.048464760119307, 1.01957599711152, 1.051874941526478, 1.010786170805605)
}

foreach ($row in $data.Keys) {
    $rowValues = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $cellRef = "{0}{1}" -f $columns[$i], $row
        $ws.Range($cellRef).Value = $rowValues[$i]
    }
}

Write-Host "applied"